$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12.87978613001893
$ws.Range("D2").Value = 6.434138702286823
$ws.Range("E2").Value = 13.10313218565756
$ws.Range("F2").Value = 47.80769781778141
$ws.Range("G2").Value = 59.65256156450099
$ws.Range("H2").Value = 21.98962784616776
$ws.Range("J2").Value = 11.03807818784017
$ws.Range("K2").Value = 20.99946033114414
$ws.Range("L2").Value = 9.132326649967723
$ws.Range("M2").Value = 21.91115033868612
$ws.Range("N2").Value = 20.94309213598616
$ws.Range("C3").Value = 12.86750037049279
$ws.Range("D3").Value = 6.428228427607915
$ws.Range("E3").Value = 13.12511409481182
$ws.Range("F3").Value = 47.87017924896716
$ws.Range("G3").Value = 59.69082686155478
$ws.Range("H3").Value = 22.04440399290321
$ws.Range("J3").Value = 11.06374489582541
$ws.Range("K3").Value = 20.75438594987395
$ws.Range("L3").Value = 9.143725407623759
$ws.Range("M3").Value = 21.83293413134121
$ws.Range("N3").Value = 21.0160937192585
$ws.Range("C4").Value = 12.86243887225138
$ws.Range("D4").Value = 6.425087244499639
$ws.Range("E4").Value = 13.14009047314802
$ws.Range("F4").Value = 47.91957307859553
$ws.Range("G4").Value = 59.73056235716435
$ws.Range("H4").Value = 22.08211580889342
$ws.Range("J4").Value = 11.08050945842339
$ws.Range("K4").Value = 20.60685778962484
$ws.Range("L4").Value = 9.151129875271762
$ws.Range("M4").Value = 21.78861388009111
$ws.Range("N4").Value = 21.06291254860909
$ws.Range("C5").Value = 12.86100215969237
$ws.Range("D5").Value = 6.42393075840102
$ws.Range("E5").Value = 13.14656585652854
$ws.Range("F5").Value = 47.94246806535997
$ws.Range("G5").Value = 59.75082703970438
$ws.Range("H5").Value = 22.09850736206935
$ws.Range("J5").Value = 11.08759441254678
$ws.Range("K5").Value = 20.54754596627496
$ws.Range("L5").Value = 9.154249567586522
$ws.Range("M5").Value = 21.77149708303594
$ws.Range("N5").Value = 21.08249521317235
$ws.Range("C6").Value = 12.86080144089408
$ws.Range("D6").Value = 6.423746214499894
$ws.Range("E6").Value = 13.14766359120934
$ws.Range("F6").Value = 47.94643663388797
$ws.Range("G6").Value = 59.75443753801191
$ws.Range("H6").Value = 22.10129093390443
$ws.Range("J6").Value = 11.08878617528474
$ws.Range("K6").Value = 20.53774791915877
$ws.Range("L6").Value = 9.154773779028449
$ws.Range("M6").Value = 21.76871221577595
$ws.Range("N6").Value = 21.08577736980606
$ws.Range("C7").Value = 12.86241695998856
$ws.Range("D7").Value = 6.425071146181055
$ws.Range("E7").Value = 13.14017629415323
$ws.Range("F7").Value = 47.91987065658213
$ws.Range("G7").Value = 59.73081918302655
$ws.Range("H7").Value = 22.08233272906886
$ws.Range("J7").Value = 11.0806039824727
$ws.Range("K7").Value = 20.60605453612496
$ws.Range("L7").Value = 9.15117153384624
$ws.Range("M7").Value = 21.78837919880781
$ws.Range("N7").Value = 21.06317460579971
$ws.Range("C8").Value = 12.87503603163437
$ws.Range("D8").Value = 6.432000169063847
$ws.Range("E8").Value = 13.11040473136593
$ws.Range("F8").Value = 47.82694821820152
$ws.Range("G8").Value = 59.66237838866449
$ws.Range("H8").Value = 22.00766686317091
$ws.Range("J8").Value = 11.04671977166647
$ws.Range("K8").Value = 20.91438655710615
$ws.Range("L8").Value = 9.136172957709929
$ws.Range("M8").Value = 21.88342021986346
$ws.Range("N8").Value = 20.96784999812488
$ws.Range("C9").Value = 12.91938278137021
$ws.Range("D9").Value = 6.449416954791685
$ws.Range("E9").Value = 13.06374433086714
$ws.Range("F9").Value = 47.73252104171043
$ws.Range("G9").Value = 59.65745193659993
$ws.Range("H9").Value = 21.8937008920056
$ws.Range("J9").Value = 10.98822461537773
$ws.Range("K9").Value = 21.53947410502304
$ws.Range("L9").Value = 9.10996425819801
$ws.Range("M9").Value = 22.09860237520207
$ws.Range("N9").Value = 20.79666872852528
$ws.Range("C10").Value = 12.96376695453056
$ws.Range("D10").Value = 6.464491427213716
$ws.Range("E10").Value = 13.03658598123773
$ws.Range("F10").Value = 47.716988857208
$ws.Range("G10").Value = 59.73308781752449
$ws.Range("H10").Value = 21.82987169071208
$ws.Range("J10").Value = 10.95006355033791
$ws.Range("K10").Value = 22.00719750999507
$ws.Range("L10").Value = 9.092641288425964
$ws.Range("M10").Value = 22.27341119834638
$ws.Range("N10").Value = 20.68038424265815
$ws.Range("C11").Value = 12.98648197011567
$ws.Range("D11").Value = 6.471831362939787
$ws.Range("E11").Value = 13.02577256609863
$ws.Range("F11").Value = 47.72166193455293
$ws.Range("G11").Value = 59.7847488334747
$ws.Range("H11").Value = 21.80517874061389
$ws.Range("J11").Value = 10.93374193423962
$ws.Range("K11").Value = 22.22096607451839
$ws.Range("L11").Value = 9.085175939164026
$ws.Range("M11").Value = 22.35636720660519
$ws.Range("N11").Value = 20.62951690493697
$ws.Range("C12").Value = 12.99544255363835
$ws.Range("D12").Value = 6.47467904304314
$ws.Range("E12").Value = 13.02189896033714
$ws.Range("F12").Value = 47.7251213932275
$ws.Range("G12").Value = 59.80679158009838
$ws.Range("H12").Value = 21.7964544325086
$ws.Range("J12").Value = 10.92771012709354
$ws.Range("K12").Value = 22.30198790313485
$ws.Range("L12").Value = 9.082408344655997
$ws.Range("M12").Value = 22.38825633495159
$ws.Range("N12").Value = 20.61054492081671
$ws.Range("C13").Value = 12.9934968390525
$ws.Range("D13").Value = 6.474062729786873
$ws.Range("E13").Value = 13.02272337932021
$ws.Range("F13").Value = 47.72430116267039
$ws.Range("G13").Value = 59.80193403670595
$ws.Range("H13").Value = 21.79830548243644
$ws.Range("J13").Value = 10.92900257207279
$ws.Range("K13").Value = 22.28453630913654
$ws.Range("L13").Value = 9.083001759837545
$ws.Range("M13").Value = 22.38136759365755
$ws.Range("N13").Value = 20.61461798945212
$ws.Range("C14").Value = 12.9872119947453
$ws.Range("D14").Value = 6.472064283201339
$ws.Range("E14").Value = 13.02544945141437
$ws.Range("F14").Value = 47.72191267245369
$ws.Range("G14").Value = 59.78651262039024
$ws.Range("H14").Value = 21.80444842304493
$ws.Range("J14").Value = 10.93324271343364
$ws.Range("K14").Value = 22.22763074408341
$ws.Range("L14").Value = 9.084947059314823
$ws.Range("M14").Value = 22.35898132282422
$ws.Range("N14").Value = 20.62795025853503
$ws.Range("C15").Value = 12.98340895690097
$ws.Range("D15").Value = 6.470849023914376
$ws.Range("E15").Value = 13.02714804379552
$ws.Range("F15").Value = 47.72066975850968
$ws.Range("G15").Value = 59.77738942377122
$ws.Range("H15").Value = 21.80829277609458
$ws.Range("J15").Value = 10.93585929232736
$ws.Range("K15").Value = 22.19278176329439
$ws.Range("L15").Value = 9.086146334466644
$ws.Range("M15").Value = 22.34533044983223
$ws.Range("N15").Value = 20.63615442339383
$ws.Range("C16").Value = 12.96233293044467
$ws.Range("D16").Value = 6.464021354315633
$ws.Range("E16").Value = 13.03732363937048
$ws.Range("F16").Value = 47.71691985882874
$ws.Range("G16").Value = 59.73005879493408
$ws.Range("H16").Value = 21.83157298908185
$ws.Range("J16").Value = 10.95115105710322
$ws.Range("K16").Value = 21.99324120821842
$ws.Range("L16").Value = 9.093137490931587
$ws.Range("M16").Value = 22.26805738442544
$ws.Range("N16").Value = 20.68374926809826
$ws.Range("C17").Value = 12.9500472385287
$ws.Range("D17").Value = 6.459955550824553
$ws.Range("E17").Value = 13.0439604585439
$ws.Range("F17").Value = 47.71762767583889
$ws.Range("G17").Value = 59.70544207511702
$ws.Range("H17").Value = 21.84696839456962
$ws.Range("J17").Value = 10.96079761581472
$ws.Range("K17").Value = 21.8710342770819
$ws.Range("L17").Value = 9.097532405823955
$ws.Range("M17").Value = 22.22151968356114
$ws.Range("N17").Value = 20.71346615988327
$ws.Range("C18").Value = 12.94321871613122
$ws.Range("D18").Value = 6.457662494245173
$ws.Range("E18").Value = 13.04792285794721
$ws.Range("F18").Value = 47.71913966303555
$ws.Range("G18").Value = 59.69290744856382
$ws.Range("H18").Value = 21.85623213299177
$ws.Range("J18").Value = 10.96644379242103
$ws.Range("K18").Value = 21.80084230668031
$ws.Range("L18").Value = 9.100099320047141
$ws.Range("M18").Value = 22.19507685903761
$ws.Range("N18").Value = 20.73074980334874
$ws.Range("C19").Value = 12.94094767116515
$ws.Range("D19").Value = 6.456893952183367
$ws.Range("E19").Value = 13.04928938726938
$ws.Range("F19").Value = 47.71984128279126
$ws.Range("G19").Value = 59.68894240177806
$ws.Range("H19").Value = 21.85943881424864
$ws.Range("J19").Value = 10.96837228957952
$ws.Range("K19").Value = 21.77709557297865
$ws.Range("L19").Value = 9.100975153482477
$ws.Range("M19").Value = 22.18618004294139
$ws.Range("N19").Value = 20.73663465108725
$ws.Range("C20").Value = 12.95133048000316
$ws.Range("D20").Value = 6.460383663299563
$ws.Range("E20").Value = 13.0432389458108
$ws.Range("F20").Value = 47.71743796092112
$ws.Range("G20").Value = 59.70789446326439
$ws.Range("H20").Value = 21.84528721419693
$ws.Range("J20").Value = 10.95976061156213
$ws.Range("K20").Value = 21.88403376808069
$ws.Range("L20").Value = 9.097060517652997
$ws.Range("M20").Value = 22.22644025507243
$ws.Range("N20").Value = 20.71028296233083
$ws.Range("C21").Value = 12.98904830010297
$ws.Range("D21").Value = 6.472649433855471
$ws.Range("E21").Value = 13.02464273795338
$ws.Range("F21").Value = 47.72256835852265
$ws.Range("G21").Value = 59.79097498510614
$ws.Range("H21").Value = 21.80262707892026
$ws.Range("J21").Value = 10.93199324523414
$ws.Range("K21").Value = 22.24434388779185
$ws.Range("L21").Value = 9.08437406875691
$ws.Range("M21").Value = 22.36554395877857
$ws.Range("N21").Value = 20.62402638071723
$ws.Range("C22").Value = 13.01578884174795
$ws.Range("D22").Value = 6.481062810157794
$ws.Range("E22").Value = 13.01377816357471
$ws.Range("F22").Value = 47.73577122841697
$ws.Range("G22").Value = 59.85972486374175
$ws.Range("H22").Value = 21.77839797061251
$ws.Range("J22").Value = 10.9147129915685
$ws.Range("K22").Value = 22.48021672139452
$ws.Range("L22").Value = 9.076428664899963
$ws.Range("M22").Value = 22.45921875837242
$ws.Range("N22").Value = 20.56934445532388
$ws.Range("C23").Value = 13.00132715544349
$ws.Range("D23").Value = 6.476536505545763
$ws.Range("E23").Value = 13.01945897418522
$ws.Range("F23").Value = 47.727823038231
$ws.Range("G23").Value = 59.82171050388184
$ws.Range("H23").Value = 21.79099479684639
$ws.Range("J23").Value = 10.9238565696873
$ws.Range("K23").Value = 22.35431455836209
$ws.Range("L23").Value = 9.080637722615849
$ws.Range("M23").Value = 22.40897624673361
$ws.Range("N23").Value = 20.59837499004975
$ws.Range("C24").Value = 12.95074959522829
$ws.Range("D24").Value = 6.460189975184599
$ws.Range("E24").Value = 13.04356468441899
$ws.Range("F24").Value = 47.71752028903887
$ws.Range("G24").Value = 59.70678069905976
$ws.Range("H24").Value = 21.84604599039407
$ws.Range("J24").Value = 10.96022912891029
$ws.Range("K24").Value = 21.87815648584431
$ws.Range("L24").Value = 9.097273732995959
$ws.Range("M24").Value = 22.22421469045913
$ws.Range("N24").Value = 20.71172146591227
$ws.Range("C25").Value = 12.90530095074928
$ws.Range("D25").Value = 6.444300959813451
$ws.Range("E25").Value = 13.0751146171039
$ws.Range("F25").Value = 47.74862974299781
$ws.Range("G25").Value = 59.64488857857626
$ws.Range("H25").Value = 21.92104484002412
$ws.Range("J25").Value = 11.00320118044778
$ws.Range("K25").Value = 21.36859348354364
$ws.Range("L25").Value = 9.116713565746098
$ws.Range("M25").Value = 22.03738976802329
$ws.Range("N25").Value = 20.84130401619938
